$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.141473094432123
$ws.Range("C2").Value = 0.105492004202171
$ws.Range("D2").Value = 0.916014941052877
$ws.Range("E2").Value = 0.213814637562741
$ws.Range("F2").Value = 0.140918641298004
$ws.Range("G2").Value = 0.000379362670713202
$ws.Range("H2").Value = 0.981557137854558
$ws.Range("I2").Value = 0.00446480681685537
$ws.Range("J2").Value = 0.988502392903
$ws.Range("K2").Value = 0.0170713201820941
$ws.Range("L2").Value = 0.122621687872067
$ws.Range("M2").Value = 0.0632952025213027
$ws.Range("N2").Value = 0.00364771798762694
$ws.Range("O2").Value = 0.179263452783938
$ws.Range("P2").Value = 0.925207190381697
$ws.Range("Q2").Value = 0.0232870316330104
$ws.Range("R2").Value = 0.968600443562507
$ws.Range("S2").Value = 0.948523403758609
$ws.Range("T2").Value = 0.233599859927629
$ws.Range("U2").Value = 0.0157581417065484
$ws.Range("V2").Value = 0.00796661608497724
$ws.Range("W2").Value = 0.0167503209991829
$ws.Range("X2").Value = 0.0213318547916423

$ws.Range("B3").Value = 0.696772499124548
$ws.Range("C3").Value = 0.348926111824443
$ws.Range("D3").Value = 0.00630325668261935
$ws.Range("E3").Value = 0.0672347379479398
$ws.Range("F3").Value = 0.699165402124431
$ws.Range("G3").Value = 0.00513598692657873
$ws.Range("H3").Value = 0.000320999182911171
$ws.Range("I3").Value = 0.00735379946305591
$ws.Range("J3").Value = 0.000320999182911171
$ws.Range("K3").Value = 0.0471285164001401
$ws.Range("L3").Value = 0.00630325668261935
$ws.Range("M3").Value = 0.0267304774133302
$ws.Range("N3").Value = 0.00671180109723357
$ws.Range("O3").Value = 0.0144157814871017
$ws.Range("P3").Value = 0.0198727675965916
$ws.Range("Q3").Value = 0.0651628341309677
$ws.Range("R3").Value = 0.0024512664876853
$ws.Range("S3").Value = 0.0317789191082059
$ws.Range("T3").Value = 0.005632076572896
$ws.Range("U3").Value = 0.0360686354616552
$ws.Range("V3").Value = 0.00265553869499241
$ws.Range("W3").Value = 0.000320999182911171
$ws.Range("X3").Value = 0.0273724757791526

$ws.Range("B4").Value = 0.136804015407961
$ws.Range("C4").Value = 0.0425469826076806
$ws.Range("D4").Value = 0.0461655188514066
$ws.Range("E4").Value = 0.687901248978639
$ws.Range("F4").Value = 0.0145033267188047
$ws.Range("G4").Value = 0.000933815804832497
$ws.Range("H4").Value = 0.00443562507295436
$ws.Range("I4").Value = 0.000204272207307109
$ws.Range("J4").Value = 0.0103011555970585
$ws.Range("K4").Value = 0.888029648651803
$ws.Range("L4").Value = 0.854616551885141
$ws.Range("M4").Value = 0.00860861445079958
$ws.Range("N4").Value = 0.00140072370724875
$ws.Range("O4").Value = 0.801739231936501
$ws.Range("P4").Value = 0.0460487918758025
$ws.Range("Q4").Value = 0.0595307575580717
$ws.Range("R4").Value = 0.0168086844869849
$ws.Range("S4").Value = 0.0157289599626474
$ws.Range("T4").Value = 0.754464806816855
$ws.Range("U4").Value = 0.013102603011556
$ws.Range("V4").Value = 0.00443562507295436
$ws.Range("W4").Value = 0.979514415781487
$ws.Range("X4").Value = 0.00525271390218279

$ws.Range("B5").Value = 0.0248920275475662
$ws.Range("C5").Value = 0.502918174390102
$ws.Range("D5").Value = 0.0314579199252947
$ws.Range("E5").Value = 0.0310201937667795
$ws.Range("F5").Value = 0.145354266370958
$ws.Range("G5").Value = 0.993521652853975
$ws.Range("H5").Value = 0.0136862378895763
$ws.Range("I5").Value = 0.987947939768881
$ws.Range("J5").Value = 0.000875452317030466
$ws.Range("K5").Value = 0.0477705147659624
$ws.Range("L5").Value = 0.0164585035601728
$ws.Range("M5").Value = 0.901336523870667
$ws.Range("N5").Value = 0.988239757207891
$ws.Range("O5").Value = 0.00449398856075639
$ws.Range("P5").Value = 0.00887125014590872
$ws.Range("Q5").Value = 0.851990194934049
$ws.Range("R5").Value = 0.0121396054628225
$ws.Range("S5").Value = 0.00396871717053811
$ws.Range("T5").Value = 0.00627407493871834
$ws.Range("U5").Value = 0.935041438076339
$ws.Range("V5").Value = 0.984942220147076
$ws.Range("W5").Value = 0.0033850822925178
$ws.Range("X5").Value = 0.94598459203922

